$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting rows 5-9 down to 6-10
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "foo"
$ws.Range("C5").Value = "0"
$ws.Range("D5").Value = "bar"
$ws.Range("E5").Value = "Baz"
